$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (col index 2)
$Bvals = @(18.22811586378697,18.07125904828559,17.97924663105714,17.94286994874786,17.93689829509151,17.97875146223086,18.17316316103315,18.58672451856243,18.90775672893032,19.05697176914917,19.11388379501154,19.10160943403612,19.06164607134108,19.03721896373105,18.89806489049617,18.81347697803073,18.76512615972327,18.74880874676969,18.82245061371383,19.0733736133093,19.23971865826845,19.15073826888531,18.8183927574385,18.47163878814503)
for ($i = 0; $i -lt $Bvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $Bvals[$i]
}

# Column C (col index 3)
$Cvals = @(6.197056111563782,6.063082965097675,5.978027535611479,5.942684148023247,5.936774795921396,5.977553618392738,6.151454133579779,6.469468434205189,6.688082070384523,6.784083823303747,6.819927446854318,6.812230788928238,6.787043002963744,6.771547908993336,6.681737504087555,6.625748224009552,6.593221041815343,6.582152775989002,6.63174198188339,6.794455221502842,6.897818258290868,6.842928568357072,6.629033256286457,6.386008921396928)
for ($i = 0; $i -lt $Cvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $Cvals[$i]
}

# Column D (col index 4)
$Dvals = @(8.632875937323508,8.641722738869484,8.647278485840962,8.649573815852971,8.649956850834821,8.647309314420632,8.635900743341073,8.614501269803004,8.599358102883734,8.592591761483728,8.590046896974343,8.590594208192268,8.592382046294043,8.593479409280132,8.599802744913907,8.603713108737194,8.605973775379553,8.60674118421808,8.603295652072985,8.59185644415756,8.584481603481549,8.588408485838409,8.603484345122201,8.620187701755629)
for ($i = 0; $i -lt $Dvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $Dvals[$i]
}

# Column F (col index 6)
$Fvals = @(42.49390039646894,42.3529784190658,42.27452276907935,42.24459811038556,42.2397532132759,42.27411088699623,42.44364629340392,42.8392972806772,43.16723786505373,43.32420334655774,43.38473173355249,43.37164793213128,43.32916140893022,43.30327807164682,43.15713387907557,43.06945174775024,43.01975387472878,43.00305408552421,43.07870984921019,43.34161142697415,43.5197648285987,43.42411209842152,43.07452204218509,42.72562048916783)
for ($i = 0; $i -lt $Fvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $Fvals[$i]
}

# Column G (col index 7)
$Gvals = @(3.739678117299573,3.74276361763592,3.744757467297474,3.74559504503969,3.745735640805469,3.744768661548654,3.740721432848829,3.733569035695276,3.728786628824841,3.726712381589509,3.725941392642641,3.726106796058088,3.726648662064826,3.726982454376524,3.72892421698855,3.730141310372184,3.730850889432067,3.731092781584428,3.730010762041392,3.726489110461469,3.724271890379621,3.72544756854559,3.730069752235354,3.735420576128568)
for ($i = 0; $i -lt $Gvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $Gvals[$i]
}

# Column K (col index 11)
$Kvals = @(14.38231300694844,14.27024283754629,14.20488810978653,14.17915078335512,14.17493191744359,14.20453734928192,14.34296987342779,14.64064215197469,14.873583012511,14.98225321454209,15.02375847768839,15.01480437975836,14.98566097875842,14.96785489643927,14.8665326258405,14.80504284198866,14.76993219602855,14.75808946130879,14.81156222137017,14.99421176738616,15.11563206645816,15.05065191808619,14.80861405772984,14.55748258158339)
for ($i = 0; $i -lt $Kvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 11).Value = $Kvals[$i]
}

# Column L (col index 12)
$Lvals = @(10.66306920895478,10.67242874732474,10.67973122093919,10.68309843793562,10.68368120600283,10.67977504742072,10.66597364402041,10.6512415297636,10.64791801145669,10.64802918043385,10.64830402943369,10.64823449528827,10.64804713173815,10.6479626566822,10.64794335357515,10.64834678771143,10.64873164614345,10.64888821653147,10.64828803018395,10.64809585371102,10.64932654740694,10.64854583805566,10.6483141180986,10.6539082050652)
for ($i = 0; $i -lt $Lvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 12).Value = $Lvals[$i]
}

# Column M (col index 13)
$Mvals = @(17.16542106310357,17.1531895000018,17.14894012836349,17.14803057748789,17.14792925344489,17.14892453072561,17.16052824229754,17.20904494051237,17.26021952816489,17.28682502319228,17.29737326164839,17.29508054278421,17.28768337721156,17.28321388432253,17.25854731037822,17.24426364475044,17.2363614359777,17.23373984974241,17.2457517663073,17.28984330020465,17.32141588015592,17.30431456304004,17.24507802162923,17.19317812743603)
for ($i = 0; $i -lt $Mvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $Mvals[$i]
}

# Column N (col index 14)
$Nvals = @(24.61716378468749,24.64478375928866,24.66355115299412,24.67165372008756,24.6730266052873,24.66365858615805,24.62631172022163,24.56742986472188,24.53292803401937,24.51913568760561,24.51418660757474,24.51524030126387,24.51872303584139,24.52089197089433,24.53386769881821,24.54231534425355,24.54735328863129,24.5490898020592,24.54139754246682,24.51769264054844,24.50379600411734,24.51106682127186,24.54181191596427,24.58182135512627)
for ($i = 0; $i -lt $Nvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 14).Value = $Nvals[$i]
}
